$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data updates
$ws.Range("B2").Value = 5.9186518994940718
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 3.6531975859042443
$ws.Range("E2").ClearContents()

# Row 3 data updates
$ws.Range("B3").Value = 4.8271078699636059
$ws.Range("C3").Value = 5.6375100864256718
$ws.Range("D3").Value = 3.1306016161220183
$ws.Range("E3").Value = 7.1032736555109457

# Update selection to reflect new selected range B1:E3
$ws.Range("B1:E3").Select()
